$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.131.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.586.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.62%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.584.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.660"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000309"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.170.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.598.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.125.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.63%  "

$ws.Range("E28").Value = "  -3.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("E34").Value = "  -2.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "572.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.01%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0806"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.397"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.52%  "

$ws.Range("E40").Value = "  +4.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.80%  "

$ws.Range("E42").Value = "  -3.40%  "

$ws.Range("E43").Value = "  -6.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.205.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0444"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.82%  "

$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.21%  "
